# Daily update at 8 AM UTC
# Appends the next day's row of data and moves the "latest row" date
# formatting (no time-of-day shown) down from the old last row to the
# new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to be the last row (row 36) reverts to the normal
# date+time number format used by every other data row.
$ws.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data in row 37.
$ws.Range("A37").Value = 45777
$ws.Range("B37").Value = 149
$ws.Range("C37").Value = 155
$ws.Range("D37").Value = 154

# The new last row gets the date-only number format.
$ws.Range("A37").NumberFormat = "YYYY-MM-DD"
